# Apply updated cryptocurrency price/volume data scraped for this run.
# Values in columns D (Price) and E (Volume 1h) are plain text cells (not numbers),
# so numeric-looking strings are written with a leading apostrophe to force Excel
# to store/keep them as text (matching the original inline-string cell type).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Bitcoin
$ws.Range('D2').Value = '27.691.27'
$ws.Range('E2').Value = '  +5.95%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '1.734.04'
$ws.Range('E3').Value = '  +4.73%  '

# Row 5: BNB
$ws.Range('D5').Value = '''227.53'
$ws.Range('E5').Value = '  +3.76%  '

# Row 6: XRP
$ws.Range('D6').Value = '''0.5454'
$ws.Range('E6').Value = '  +3.67%  '

# Row 8: Cardano
$ws.Range('D8').Value = '''0.2743'
$ws.Range('E8').Value = '  +2.00%  '

# Row 9: Dogecoin
$ws.Range('D9').Value = '''0.06714'
$ws.Range('E9').Value = '  +5.28%  '

# Row 10: Solana
$ws.Range('D10').Value = '''21.84'
$ws.Range('E10').Value = '  +6.14%  '

# Row 11: TRON
$ws.Range('D11').Value = '''0.07774'
$ws.Range('E11').Value = '  +1.08%  '

# Row 12: Polkadot
$ws.Range('D12').Value = '''4.697'
$ws.Range('E12').Value = '  +1.81%  '

# Row 13: WrappedliquidstakedEther2.0
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.747.29'
$ws.Range('E13').Value = '  +6.34%  '

# Row 14: WrappedEther
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').Value = '1.973.95'
$ws.Range('E14').Value = '  +4.79%  '

# Row 15: Polygon
$ws.Range('D15').Value = '''0.5983'
$ws.Range('E15').Value = '  +6.15%  '

# Row 16: ShibaInu
$ws.Range('D16').Value = '0.0₅8422'
$ws.Range('E16').Value = '  +1.96%  '

# Row 17: Litecoin
$ws.Range('D17').Value = '''69.12'
$ws.Range('E17').Value = '  +5.08%  '

# Row 18: WrappedBTC
$ws.Range('D18').Value = '27.709.24'
$ws.Range('E18').Value = '  +6.05%  '

# Row 19: BitcoinCash
$ws.Range('D19').Value = '''226.53'
$ws.Range('E19').Value = '  +18.72%  '

# Row 20: Uniswap
$ws.Range('D20').Value = '''4.828'
$ws.Range('E20').Value = '  +2.90%  '

# Row 21: Dai
$ws.Range('D21').Value = '''1.003'
$ws.Range('E21').Value = '  -0.08%  '

# Row 22: Avalanche
$ws.Range('D22').Value = '''10.89'
$ws.Range('E22').Value = '  +5.06%  '

# Row 23: Chainlink
$ws.Range('D23').Value = '''6.222'
$ws.Range('E23').Value = '  +3.64%  '

# Row 24: BinanceUSD
$ws.Range('E24').Value = '  -0.14%  '

# Row 25: Monero
$ws.Range('D25').Value = '''148.15'
$ws.Range('E25').Value = '  -2.07%  '

# Row 26: Toncoin
$ws.Range('E26').Value = '  +13.35%  '

# Row 28: Cosmos
$ws.Range('D28').Value = '''7.471'
$ws.Range('E28').Value = '  +2.65%  '

# Row 29: EthereumClassic
$ws.Range('D29').Value = '''17.10'
$ws.Range('E29').Value = '  +6.65%  '

# Row 30: Hedera
$ws.Range('D30').Value = '''0.05707'
$ws.Range('E30').Value = '  +1.00%  '

# Row 31: PancakeSwap
$ws.Range('D31').Value = '''1.312'
$ws.Range('E31').Value = '  +2.68%  '

# Row 32: InternetComputer(DFINITY)
$ws.Range('D32').Value = '''3.703'
$ws.Range('E32').Value = '  +5.93%  '

# Row 33: Filecoin
$ws.Range('E33').Value = '  +3.96%  '

# Row 34: LidoDAOToken
$ws.Range('D34').Value = '''1.687'
$ws.Range('E34').Value = '  +6.50%  '

# Row 35: ARBITRUM
$ws.Range('D35').Value = '''0.9755'
$ws.Range('E35').Value = '  +2.74%  '

# Row 36: MXToken
$ws.Range('E36').Value = '  +1.96%  '

# Row 37: HuobiToken
$ws.Range('D37').Value = '''2.437'
$ws.Range('E37').Value = '  +1.14%  '

# Row 38: ImmutableX
$ws.Range('D38').Value = '''0.5988'
$ws.Range('E38').Value = '  +3.50%  '

# Row 39: VeChain
$ws.Range('D39').Value = '''0.01670'
$ws.Range('E39').Value = '  +4.26%  '

# Row 40: FraxShare
$ws.Range('D40').Value = '''5.934'
$ws.Range('E40').Value = '  -0.70%  '

# Row 41: TrustWalletToken
$ws.Range('D41').Value = '''0.8505'
$ws.Range('E41').Value = '  +1.94%  '

# Row 42: Maker
$ws.Range('D42').Value = '1.049.09'
$ws.Range('E42').Value = '  +2.36%  '

# Row 43: PaxDollar
$ws.Range('D43').Value = '''1.003'
$ws.Range('E43').Value = '  -0.02%  '

# Row 44: Quant
$ws.Range('D44').Value = '''101.66'
$ws.Range('E44').Value = '  +0.18%  '

# Row 45: RocketPoolETH
$ws.Range('D45').Value = '1.878.33'

# Row 46: BabyDogeCoin
$ws.Range('E46').Value = '  +10.90%  '

# Row 47: Aave
$ws.Range('D47').Value = '''59.70'
$ws.Range('E47').Value = '  +2.11%  '

# Row 48: EnergySwap
$ws.Range('D48').Value = '''8.287'
$ws.Range('E48').Value = '  +2.77%  '

# Row 49: Mantle
$ws.Range('D49').Value = '''0.4428'
$ws.Range('E49').Value = '  +2.01%  '

# Row 50: Cronos
$ws.Range('E50').Value = '  -0.21%  '

# Row 51: Frax
$ws.Range('D51').Value = '''0.9982'
$ws.Range('E51').Value = '  -0.56%  '
